# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted into the "Cebollín" data set
# for Feria Lagunitas de Puerto Montt. The new record is placed right
# before the existing row 310, pushing that row (and every row after it)
# down by one. All the other data for the new row is identical to the
# row that used to be at position 310 (same category/quality/prices/
# origin/etc.), except for the date (column D) and the volume
# (column J), which carry the new observation's own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 310; this shifts the former rows
# 310..411 down to 311..412 and copies formatting from the row above,
# matching Excel's normal "Insert" behaviour.
$ws.Rows(310).Insert()

# The row that used to be 310 now lives at 311. Duplicate it into the
# freshly inserted row 310 so every column starts out identical.
$ws.Range("A311:R311").Copy($ws.Range("A310:R310"))

# Overwrite the two fields that differ for this new observation.
$ws.Cells.Item(310, 4).Value2 = 44985   # D310 - Fecha
$ws.Cells.Item(310, 10).Value2 = 180    # J310 - Volumen
